# Append a duplicate copy of the leading data rows to the end of each
# sheet's used range (CustomerSummary gains rows 42-61 = a copy of rows
# 2-21; TransactionSummary gains rows 78-96 = a copy of rows 2-20).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("CustomerSummary")
$ws1Rows = $ws1.UsedRange.Rows.Count
$ws1Src = $ws1.Range("A2:F21")
$ws1Src.Copy()
$ws1Dest = $ws1.Cells.Item($ws1Rows + 1, 1)
$ws1.Paste($ws1Dest)

$ws2 = $wb.Worksheets.Item("TransactionSummary")
$ws2Rows = $ws2.UsedRange.Rows.Count
$ws2Src = $ws2.Range("A2:F20")
$ws2Src.Copy()
$ws2Dest = $ws2.Cells.Item($ws2Rows + 1, 1)
$ws2.Paste($ws2Dest)
